$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet tab to reflect the new "through" date
$ws.Name = "Through 2022-02-24"

# Update the header label in I1 (shared string "2022 (through 02-23)" -> "2022 (through 02-24)")
$ws.Range("I1").Value = "2022 (through 02-24)"

# Update the data values for the new day of data (2022-03-04 commit)
$ws.Range("I3").Value = 121
$ws.Range("I14").Value = 280
